$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string as TEXT (t="s") without touching any
# cell style, by round-tripping it through a scratch cell + copy/paste-values.
# (A plain ".Value = '977448714'" gets auto-coerced to a number by Excel.)
function Set-TextValue($cellRef, $text) {
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '=TEXT(' + $text + ',"0")'
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# Row 3 (existing row): fill in the previously blank Response cell.
$ws.Range("G3").Value = "Something something"

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "22962f26-a8ad-4c36-b5eb-9449d90c0ed4"
$ws.Range("C4").Value = "Mahesh"
Set-TextValue "D4" 977448714
$ws.Range("E4").Value = "maheshanna@gmail.com"
$ws.Range("F4").Value = "How are you ?"
$ws.Range("G4").Value = "I am good!!"

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "12962f26-a8ad-4c36-b5eb-9449d90c0ed4"
$ws.Range("C5").Value = "Jamanalal"
Set-TextValue "D5" 977448724
$ws.Range("E5").Value = "maheshanna@gmail.com"
$ws.Range("F5").Value = "How are you ?"
$ws.Range("G5").Value = "What are the data"

# Row 6
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "54ebd306-48fd-49f9-b7c5-ffb134f0dd77"
$ws.Range("C6").Value = "fadsf"
$ws.Range("D6").Value = "adsfa"
$ws.Range("E6").Value = "fasdfa"
$ws.Range("F6").Value = "afasfaf"
$ws.Range("G6").Value = ""

# Row 7
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "4b4ce39b-75bd-4afd-839d-2ac907064aa8"
$ws.Range("C7").Value = "Biraj Dey"
Set-TextValue "D7" 8789445445
$ws.Range("E7").Value = "biraj123@gmail.com"
$ws.Range("F7").Value = "sdfkajflksd"
$ws.Range("G7").Value = ""

$ws.Range("K12").Select()
